# "Refined metadata to be additional tab"
#
# 1. Re-stamp the panel-query timestamps in the "data" sheet's F column
#    (rows 2-43) with their new values.
# 2. Add a new "metadata" worksheet (after "data") summarising the panel
#    query itself (name/id/version/etc.).

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Update data!F2:F43 timestamps -------------------------------------
$newTimes = @(
  "2021-10-05 14:35:48.507450",
  "2021-10-05 14:35:48.507458",
  "2021-10-05 14:35:48.507461",
  "2021-10-05 14:35:48.507464",
  "2021-10-05 14:35:48.507467",
  "2021-10-05 14:35:48.507469",
  "2021-10-05 14:35:48.507472",
  "2021-10-05 14:35:48.507474",
  "2021-10-05 14:35:48.507477",
  "2021-10-05 14:35:48.507480",
  "2021-10-05 14:35:48.507482",
  "2021-10-05 14:35:48.507485",
  "2021-10-05 14:35:48.507488",
  "2021-10-05 14:35:48.507491",
  "2021-10-05 14:35:48.507493",
  "2021-10-05 14:35:48.507496",
  "2021-10-05 14:35:48.507499",
  "2021-10-05 14:35:48.507501",
  "2021-10-05 14:35:48.507504",
  "2021-10-05 14:35:48.507507",
  "2021-10-05 14:35:48.507510",
  "2021-10-05 14:35:48.507512",
  "2021-10-05 14:35:48.507515",
  "2021-10-05 14:35:48.507517",
  "2021-10-05 14:35:48.507520",
  "2021-10-05 14:35:48.507523",
  "2021-10-05 14:35:48.507526",
  "2021-10-05 14:35:48.507528",
  "2021-10-05 14:35:48.507531",
  "2021-10-05 14:35:48.507533",
  "2021-10-05 14:35:48.507536",
  "2021-10-05 14:35:48.507538",
  "2021-10-05 14:35:48.507541",
  "2021-10-05 14:35:48.507544",
  "2021-10-05 14:35:48.507547",
  "2021-10-05 14:35:48.507549",
  "2021-10-05 14:35:48.507552",
  "2021-10-05 14:35:48.507554",
  "2021-10-05 14:35:48.507557",
  "2021-10-05 14:35:48.507559",
  "2021-10-05 14:35:48.507562",
  "2021-10-05 14:35:48.507565"
)

for ($i = 0; $i -lt $newTimes.Count; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- 2. Add the "metadata" worksheet, right after "data" ------------------
$ws2 = $wb.Worksheets.Add($null, $dataSheet)
$ws2.Name = "metadata"

# Match the header-row look-and-feel (bold / centered / bordered) used on
# the "data" sheet by copying its cell formatting across.
$dataSheet.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("B1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Susceptibility to Viral Infections"
$ws2.Range("C2").Value = 237

# These look numeric/date-like but must be stored as plain text, matching
# the source data export. Forcing a text number-format before assignment
# keeps them as strings; resetting the style afterwards avoids leaving a
# stray "Text" number-format behind on the cell.
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "0.77"
$ws2.Range("D2").Style = "Normal"

$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "2021-07-25T07:25:31.259899Z"
$ws2.Range("E2").Style = "Normal"

$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "2021-10-05 14:35:48.503746"
$ws2.Range("F2").Style = "Normal"

$ws2.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/237/?format=json"

Write-Output "metadata tab added; data!F2:F43 timestamps refreshed"
